$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to retain text formatting so values such as
# "37.102.29" or "2.41" are not auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Updated prices (column D)
$ws.Range("D2").Value = "37.102.29"
$ws.Range("D3").Value = "2.029.03"
$ws.Range("D5").Value = "226.36"
$ws.Range("D8").Value = "55.36"
$ws.Range("D9").Value = "0.382"
$ws.Range("D12").Value = "2.329.30"
$ws.Range("D13").Value = "14.34"
$ws.Range("D14").Value = "20.57"
$ws.Range("D15").Value = "0.747"
$ws.Range("D16").Value = "5.17"
$ws.Range("D17").Value = "2.031.23"
$ws.Range("D18").Value = "37.066.73"
$ws.Range("D19").Value = "6.13"
$ws.Range("D20").Value = "68.81"
$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("D22").Value = "226.55"
$ws.Range("D24").Value = "2.41"
$ws.Range("D25").Value = "2.28"
$ws.Range("D26").Value = "9.42"
$ws.Range("D27").Value = "167.61"
$ws.Range("D29").Value = "18.81"
$ws.Range("D31").Value = "0.117"
$ws.Range("D34").Value = "4.46"
$ws.Range("D39").Value = "5.44"
$ws.Range("D40").Value = "1.505.59"
$ws.Range("D42").Value = "16.74"
$ws.Range("D43").Value = "2.84"
$ws.Range("D44").Value = "95.49"
$ws.Range("D45").Value = "0.0928"
$ws.Range("D46").Value = "1.13"
$ws.Range("D47").Value = "7.22"
$ws.Range("D49").Value = "2.93"
$ws.Range("D50").Value = "3.70"
$ws.Range("D51").Value = "2.215.62"

# Restore the default (General) style on column D now that the text values
# are safely stored, so no stray number-format style remains on the cells.
$priceRange.Style = "Normal"

# Updated 1-hour volume percentages (column E) -- these already contain
# non-numeric characters (%, padding spaces) so Excel keeps them as text.
$ws.Range("E2").Value = "  -1.56%  "
$ws.Range("E3").Value = "  -2.59%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("E5").Value = "  -3.06%  "
$ws.Range("E6").Value = "  -4.25%  "
$ws.Range("E8").Value = "  -4.27%  "
$ws.Range("E9").Value = "  -2.19%  "
$ws.Range("E10").Value = "  +2.42%  "
$ws.Range("E11").Value = "  -3.44%  "
$ws.Range("E12").Value = "  -2.59%  "
$ws.Range("E13").Value = "  -5.51%  "
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("E15").Value = "  -3.44%  "
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("E18").Value = "  -1.66%  "
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("E25").Value = "  -4.36%  "
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("E28").Value = "  -4.69%  "
$ws.Range("E29").Value = "  -3.52%  "
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("E31").Value = "  -4.73%  "
$ws.Range("E32").Value = "  -3.33%  "
$ws.Range("E33").Value = "  -3.58%  "
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  +2.97%  "
$ws.Range("E41").Value = "  -6.86%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("E44").Value = "  -5.40%  "
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("E46").Value = "  -5.68%  "
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("E48").Value = "  -3.97%  "
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("E50").Value = "  -8.80%  "
$ws.Range("E51").Value = "  -2.65%  "
